$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010-18")

# New row 9 data (CW3M Baseline 2010-18 C379 run)
$ws.Cells.Item(9, 1).Value = "CW3M"
$ws.Cells.Item(9, 2).Value = "Baseline 2010-18 C379"
$ws.Cells.Item(9, 3).Value = "2010-18"

$ws.Cells.Item(9, 4).Value = 687.27896466666664
$ws.Cells.Item(9, 5).Value = 2094.2995878888887
$ws.Cells.Item(9, 6).Value = 5.8066811111111116
$ws.Cells.Item(9, 7).Value = 195.47808666666668
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 5.9917683333333338
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 533.31639266666662
$ws.Cells.Item(9, 12).Value = 89.974237444444441
$ws.Cells.Item(9, 13).Value = 1672.4987521111111
$ws.Cells.Item(9, 14).Value = 693.36477322222231
$ws.Cells.Item(9, 15).Value = 15597.417643111112
$ws.Cells.Item(9, 16).Value = 2216.8192002222222
$ws.Cells.Item(9, 17).Value = 0.29906677777777779
$ws.Cells.Item(9, 18).Value = 0.00000077777777777776053

# Match number formats used by the rest of the table
$ws.Range("D9:N9").NumberFormat = "0.00"
$ws.Range("O9:P9").NumberFormat = "0"
$ws.Range("Q9").NumberFormat = "0.00"
$ws.Range("R9").NumberFormat = "0.000000"

# Move the active selection down to the next empty row, as in the saved file
$ws.Range("B10").Select()
